$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.705.68'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  -0.39%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.637.30'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  -0.83%  '

$ws.Range("E4").Value = '  -0.28%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.54'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +0.41%  '

$ws.Range("E6").Value = '  -1.08%  '

$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("E8").Value = '  -0.51%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0622'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  -0.87%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.07'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -0.77%  '

$ws.Range("E11").Value = '  +0.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.864.59'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  -0.87%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.636.45'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -0.67%  '

$ws.Range("E14").Value = '  -1.39%  '

$ws.Range("E15").Value = '  -1.58%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.44'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  -1.71%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.691.81'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -0.47%  '

$ws.Range("E18").Value = '  -2.54%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '211.26'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  -3.59%  '

$ws.Range("E20").Value = '  -0.19%  '

$ws.Range("E21").Value = '  -0.97%  '

$ws.Range("E22").Value = '  -1.44%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.30'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  -2.63%  '

$ws.Range("E24").Value = '  -2.70%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.65'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  -0.18%  '

$ws.Range("E26").Value = '  -0.28%  '

$ws.Range("E27").Value = '  -2.00%  '

$ws.Range("E28").Value = '  -0.74%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.54'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  -1.38%  '

$ws.Range("E30").Value = '  -2.76%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.19'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  +0.40%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.34'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  -0.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.98'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  -1.33%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.271.79'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  -0.99%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.52'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  -1.38%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.43'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  -0.28%  '

$ws.Range("E37").Value = '  -2.20%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.528'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  -2.02%  '

$ws.Range("E39").Value = '  -3.02%  '

$ws.Range("E40").Value = '  -0.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.802'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  -1.81%  '

$ws.Range("E42").Value = '  -2.75%  '

$ws.Range("E43").Value = '  -3.78%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.774.91'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  -0.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.45'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -0.69%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.25'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +0.49%  '

$ws.Range("E47").Value = '  -1.53%  '

$ws.Range("E48").Value = '  +0.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.53'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  -2.98%  '

$ws.Range("E50").Value = '  -1.05%  '

$ws.Range("E51").Value = '  -0.56%  '
